$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: new "LED current-limiting resistor" calculation block ---

# Row 1 headers (order chosen so new shared-string indices land as
# Current=119, Voltage Drop=120, Resistance=121, Voltage Supply=122)
$ws2.Range("A7").Value = "Current"
$ws2.Range("C1").Value = "Voltage Drop"
$ws2.Range("D1").Value = "Resistance"
$ws2.Range("B1").Value = "Voltage Supply"

# B7 (Current, amps) gets the scientific-style "0.00E+00" number format
# (built-in numFmtId 11) - set before the D column format so it becomes
# cellXfs index 11.
$ws2.Range("B7").Value = 0.01
$ws2.Range("B7").NumberFormat = "0.00E+00"

# Data rows: Red, Green, Yellow, Blue LEDs
$ws2.Range("A2").Value = "Red"
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 1.85

$ws2.Range("A3").Value = "Green"
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = 2.2

$ws2.Range("A4").Value = "Yellow"
$ws2.Range("B4").Value = 5
$ws2.Range("C4").Value = 2

$ws2.Range("A5").Value = "Blue"
$ws2.Range("B5").Value = 5
$ws2.Range("C5").Value = 3.3

# Resistance formulas (D2:D5) - set after B7's format so this run of
# cells becomes cellXfs index 12 (numFmtId 2, "0.00").
$ws2.Range("D2").Formula = '=(B2-C2)/$B$7'
$ws2.Range("D3").Formula = '=(B3-C3)/$B$7'
$ws2.Range("D4").Formula = '=(B4-C4)/$B$7'
$ws2.Range("D5").Formula = '=(B5-C5)/$B$7'
$ws2.Range("D2:D5").NumberFormat = "0.00"

# Column widths for C and D
$ws2.Columns.Item(3).ColumnWidth = 12.8
$ws2.Columns.Item(4).ColumnWidth = 15.8

# --- View / selection changes ---
# Sheet1 scrolls back to the top-left and the selection moves off the
# bottom of the BOM to G7; Sheet2 becomes the active tab with E9 selected.
$ws1.Range("G7").Select() | Out-Null
$ws2.Range("E9").Select() | Out-Null
$ws2.Activate()
